$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$wsVentasGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasGrupo.Range("L17").Value = 177.19

# Sheet: VENTA MENSUAL
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F17").Value = 306.85
$wsVentaMensual.Range("F30").Value = 1687.38

# Sheet: CUMPLIMIENTO MENSUAL
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D16").Value = 1225.24
$wsCumplimiento.Range("E16").Value = 17573.37
$wsCumplimiento.Range("F16").Value = 0.06517715937508145
$wsCumplimiento.Range("D19").Value = 1681.62
$wsCumplimiento.Range("E19").Value = 27856.17107555788
$wsCumplimiento.Range("F19").Value = 0.05693113597081122
